$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: this used to be the last data row of its group; it becomes a
# "group separator" row (matching the look of rows 4 and 7) now that two more
# rows follow it. Copy the formatting (borders/fonts, style 6/7) from row 7,
# which already has that exact look, without touching the values in row 10.
$ws.Range("A7:E7").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# --- New rows 11 & 12 (row 11 plain like row 9, row 12 taller/wrapped like row 6) ---
$ws.Range("B11").Value = 325
$ws.Range("B12").Value = 328

$ws.Range("C11").Value = " Ah, welcome back! ♪"
$ws.Range("C12").Value = " It sounds very much like it was\na fine expedition for your guild! ♪"

$ws.Range("D11").Value = " Ах, с возвращением! ♪"
$ws.Range("D12").Value = " Похоже, ваша экспедиция\nувенчалась успехом! ♪"

$ws.Range("E11").Value = " Àö, ò âïèâñàþåîéåí! ♪"
$ws.Range("E12").Value = " Ðïöïçå, âàšà üëòðåäéøéÿ\nôâåîœàìàòû ôòðåöïí! ♪"

$ws.Rows.Item(12).RowHeight = 21.6

# --- Update selection to follow the newly active cell, like the original edit ---
[void]$ws.Range("D11").Select()
